$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 90152
$ws.Range("J3").Value = 90152
$ws.Range("L3").Value = 90152
$ws.Range("N3").Value = -90380
$ws.Range("H18").Value = 342
$ws.Range("I18").Value = 342
$ws.Range("K18").Value = 342
$ws.Range("M18").Value = -58
$ws.Range("H74").Value = 5492.5264
$ws.Range("J74").Value = 5990.909
$ws.Range("L74").Value = 5990.909
$ws.Range("N74").Value = -7862.909
$ws.Range("H77").Value = 5492.5264
$ws.Range("J77").Value = 5990.909
$ws.Range("L77").Value = 29954.545
$ws.Range("N77").Value = -39314.545
$ws.Range("H102").Value = 90152
$ws.Range("J102").Value = 90152
$ws.Range("L102").Value = 90152
$ws.Range("N102").Value = -96642

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 1099.2727
$ws.Range("I15").Value = 1099.2727
$ws.Range("K15").Value = 1099.2727
$ws.Range("M15").Value = -749.2727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 69451496
$ws.Range("I31").Value = 4935.5
$ws.Range("J31").Value = 125008750
$ws.Range("K31").Value = 4935.5
$ws.Range("L31").Value = 125008750
$ws.Range("M31").Value = -4640.5
$ws.Range("N31").Value = -125009340
$ws.Range("H34").Value = 69451496
$ws.Range("I34").Value = 4935.5
$ws.Range("J34").Value = 125008750
$ws.Range("K34").Value = 4935.5
$ws.Range("L34").Value = 125008750
$ws.Range("M34").Value = -4733.5
$ws.Range("N34").Value = -125009154
$ws.Range("H62").Value = 3204.5881
$ws.Range("I62").Value = 3237.6875
$ws.Range("K62").Value = 3237.6875
$ws.Range("M62").Value = -2613.6875
$ws.Range("H65").Value = 3204.5881
$ws.Range("I65").Value = 3237.6875
$ws.Range("K65").Value = 16188.4375
$ws.Range("M65").Value = -13068.4375
$ws.Range("H99").Value = 3980.5
$ws.Range("I99").Value = 3478.818
$ws.Range("J99").Value = 9499
$ws.Range("K99").Value = 3478.818
$ws.Range("L99").Value = 9499
$ws.Range("M99").Value = -1980.818
$ws.Range("N99").Value = -12495
$ws.Range("H126").Value = 3980.5
$ws.Range("I126").Value = 3478.818
$ws.Range("J126").Value = 9499
$ws.Range("K126").Value = 10436.454
$ws.Range("L126").Value = 28497
$ws.Range("M126").Value = -7966.454000000002
$ws.Range("N126").Value = -33437
$ws.Range("H131").Value = 22464
$ws.Range("J131").Value = 22464
$ws.Range("L131").Value = 22464
$ws.Range("N131").Value = -32544
$ws.Range("H138").Value = 200000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 55495.668
$ws.Range("J140").Value = 67889
$ws.Range("L140").Value = 67889
$ws.Range("N140").Value = -78249
$ws.Range("H141").Value = 94045
$ws.Range("J141").Value = 94045
$ws.Range("L141").Value = 94045
$ws.Range("N141").Value = -104405

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 698.7059
$ws.Range("I23").Value = 402.9
$ws.Range("J23").Value = 1121.2858
$ws.Range("K23").Value = 1208.7
$ws.Range("L23").Value = 3363.8574
$ws.Range("M23").Value = -973.6999999999998
$ws.Range("N23").Value = -3833.8574
$ws.Range("H34").Value = 759.95
$ws.Range("J34").Value = 1299.9
$ws.Range("L34").Value = 3899.7
$ws.Range("N34").Value = -4067.7
$ws.Range("H128").Value = 185993.33
$ws.Range("I128").Value = 185993.33
$ws.Range("K128").Value = 557979.99
$ws.Range("M128").Value = -552999.99

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 19981.666
$ws.Range("I20").Value = 18000
$ws.Range("K20").Value = 18000
$ws.Range("M20").Value = -17755
$ws.Range("H24").Value = 19332
$ws.Range("I24").Value = 19332
$ws.Range("K24").Value = 19332
$ws.Range("M24").Value = -19159
$ws.Range("H122").Value = 45456610
$ws.Range("I122").Value = 2206.2856
$ws.Range("K122").Value = 6618.8568
$ws.Range("M122").Value = -4168.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3133.7932
$ws.Range("I7").Value = 3038.1904
$ws.Range("J7").Value = 3384.75
$ws.Range("K7").Value = 3038.1904
$ws.Range("L7").Value = 3384.75
$ws.Range("M7").Value = -2926.1904
$ws.Range("N7").Value = -3608.75
$ws.Range("H20").Value = 4581.25
$ws.Range("I20").Value = 4100
$ws.Range("K20").Value = 4100
$ws.Range("M20").Value = -3874
$ws.Range("H22").Value = 2397.7334
$ws.Range("J22").Value = 2469
$ws.Range("L22").Value = 2469
$ws.Range("N22").Value = -3059
$ws.Range("H27").Value = 2397.7334
$ws.Range("J27").Value = 2469
$ws.Range("L27").Value = 2469
$ws.Range("N27").Value = -2683
$ws.Range("H40").Value = 2823.6667
$ws.Range("J40").Value = 5798.2
$ws.Range("L40").Value = 5798.2
$ws.Range("N40").Value = -6070.2
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -35980
$ws.Range("H52").Value = 35000
$ws.Range("J52").Value = 35000
$ws.Range("L52").Value = 35000
$ws.Range("N52").Value = -35466
$ws.Range("H126").Value = 3133.7932
$ws.Range("I126").Value = 3038.1904
$ws.Range("J126").Value = 3384.75
$ws.Range("K126").Value = 9114.5712
$ws.Range("L126").Value = 10154.25
$ws.Range("M126").Value = -6644.5712
$ws.Range("N126").Value = -15094.25
$ws.Range("H132").Value = 105272000
$ws.Range("I132").Value = 2837.182
$ws.Range("K132").Value = 8511.545999999998
$ws.Range("M132").Value = -5981.545999999998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H96").Value = 6555.5713
$ws.Range("I96").Value = 5378
$ws.Range("K96").Value = 5378
$ws.Range("M96").Value = -4005
$ws.Range("H132").Value = 1243.4706
$ws.Range("I132").Value = 1275.2667
$ws.Range("K132").Value = 3825.800099999999
$ws.Range("M132").Value = -1295.800099999999
